$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44181
$ws.Range("J2").Value = 12000
$ws.Range("D3").Value = 44204
$ws.Range("D4").Value = 44186
$ws.Range("J4").Value = 10000
$ws.Range("D5").Value = 44189
$ws.Range("J5").Value = 16000
$ws.Range("D6").Value = 44159
$ws.Range("J6").Value = 7000
$ws.Range("D7").Value = 44232
$ws.Range("J7").Value = 16000
$ws.Range("D8").Value = 44231
$ws.Range("J8").Value = 12000
$ws.Range("D9").Value = 44188
$ws.Range("J9").Value = 12000
$ws.Range("D10").Value = 44230
$ws.Range("D11").Value = 44229
$ws.Range("J11").Value = 16000
$ws.Range("D12").Value = 44245
$ws.Range("J12").Value = 9000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("D13").Value = 44245
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 25
$ws.Range("D14").Value = 44161
$ws.Range("D15").Value = 44214
$ws.Range("J15").Value = 7000
$ws.Range("O15").Value = "Provincia de Chacabuco"
$ws.Range("D16").Value = 44167
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("O16").Value = "Provincia de Chacabuco"
$ws.Range("P16").Value = 30
$ws.Range("D17").Value = 44210
$ws.Range("J17").Value = 8800
$ws.Range("K17").Value = 2500
$ws.Range("M17").Value = 2750
$ws.Range("P17").Value = 28
$ws.Range("D18").Value = 44187
$ws.Range("J18").Value = 12000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("P18").Value = 30
$ws.Range("D19").Value = 44166
$ws.Range("J19").Value = 7000
$ws.Range("D21").Value = 44209
$ws.Range("J21").Value = 7000
$ws.Range("K21").Value = 2500
$ws.Range("M21").Value = 2750
$ws.Range("P21").Value = 28
$ws.Range("D22").Value = 44160
$ws.Range("J22").Value = 7000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("P22").Value = 30
$ws.Range("D23").Value = 44215
$ws.Range("J23").Value = 16000
$ws.Range("D24").Value = 44162
